# Add a new "addItemsToCartTest" test case to both the RUNMANAGER and DATA
# sheets, and refresh the DATA sheet's old spare row (previously a
# locked_out_user data row) so it now backs the new add-to-cart test.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# --- DATA sheet: replace the old (unused) row 3 with the new test's data row.
# Delete it first so the freshly written row 3 starts out with default
# (unstyled) formatting rather than inheriting the old row's styles.
$ws2.Rows.Item(3).Delete()
$ws2.Range("A3").Value = "addItemsToCartTest"
$ws2.Range("B3").Value = "yes"
$ws2.Range("C3").Value = "chrome"
$ws2.Range("D3").Value = "standard_user"
$ws2.Range("E3").Value = "secret_sauce"

# --- RUNMANAGER sheet: append the new test case as row 3.
$ws1.Range("A3").Value = "addItemsToCartTest"
$ws1.Range("B3").Value = "To check whether the user can add items to cart"
$ws1.Range("C3").Value = "yes"
$ws1.Range("D3").Value = "1"
$ws1.Range("E3").Value = "1"
$ws1.Range("A3").Style = "Normal"

# --- View/selection state: DATA becomes the active/selected tab, RUNMANAGER
# loses tabSelected. Selections move too.
$ws1.Activate()
$ws1.Range("D6").Select()
$ws2.Activate()
$ws2.Range("A3").Select()

# --- Page setup on DATA sheet.
$ws2.PageSetup.Orientation = "xlPortrait"
